$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.707.79'
$ws.Range("E2").Value = '  -1.08%  '

# Row 3
$ws.Range("D3").Value = '1.741.31'
$ws.Range("E3").Value = '  -2.06%  '

# Row 4
$ws.Range("E4").Value = '  +0.81%  '

# Row 5
$ws.Range("D5").Value = "'330.90"
$ws.Range("E5").Value = '  -0.65%  '

# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = '  +0.95%  '

# Row 7
$ws.Range("D7").Value = "'0.3847"
$ws.Range("E7").Value = '  +0.80%  '

# Row 8
$ws.Range("D8").Value = "'0.3340"
$ws.Range("E8").Value = '  -2.40%  '

# Row 9
$ws.Range("D9").Value = "'45.34"
$ws.Range("E9").Value = '  -4.95%  '

# Row 10
$ws.Range("D10").Value = "'1.100"
$ws.Range("E10").Value = '  -4.02%  '

# Row 11
$ws.Range("D11").Value = "'0.07170"
$ws.Range("E11").Value = '  -3.60%  '

# Row 12
$ws.Range("E12").Value = '  +0.82%  '

# Row 13
$ws.Range("D13").Value = "'22.14"
$ws.Range("E13").Value = '  -3.14%  '

# Row 14
$ws.Range("D14").Value = "'6.093"
$ws.Range("E14").Value = '  -4.65%  '

# Row 15
$ws.Range("D15").Value = '1.741.35'
$ws.Range("E15").Value = '  -1.99%  '

# Row 17
$ws.Range("E17").Value = '  -2.98%  '

# Row 18
$ws.Range("D18").Value = "'0.06584"
$ws.Range("E18").Value = '  -0.70%  '

# Row 19
$ws.Range("E19").Value = '  +0.63%  '

# Row 20
$ws.Range("D20").Value = "'78.21"
$ws.Range("E20").Value = '  -5.51%  '

# Row 21
$ws.Range("D21").Value = "'16.59"
$ws.Range("E21").Value = '  -5.17%  '

# Row 22
$ws.Range("D22").Value = "'6.145"
$ws.Range("E22").Value = '  -4.47%  '

# Row 23
$ws.Range("D23").Value = '27.709.95'
$ws.Range("E23").Value = '  -1.20%  '

# Row 24
$ws.Range("D24").Value = "'11.50"
$ws.Range("E24").Value = '  -5.12%  '

# Row 25
$ws.Range("D25").Value = "'2.393"
$ws.Range("E25").Value = '  +0.43%  '

# Row 26
$ws.Range("D26").Value = "'155.14"
$ws.Range("E26").Value = '  +0.73%  '

# Row 27
$ws.Range("D27").Value = "'19.62"
$ws.Range("E27").Value = '  -6.00%  '

# Row 28
$ws.Range("D28").Value = "'2.252"
$ws.Range("E28").Value = '  -7.97%  '

# Row 29
$ws.Range("D29").Value = '1.939.50'
$ws.Range("E29").Value = '  -2.01%  '

# Row 30
$ws.Range("D30").Value = "'1.256"
$ws.Range("E30").Value = '  -12.90%  '

# Row 31
$ws.Range("D31").Value = "'128.10"
$ws.Range("E31").Value = '  -5.16%  '

# Row 32
$ws.Range("D32").Value = "'4.015"
$ws.Range("E32").Value = '  +1.58%  '

# Row 33
$ws.Range("D33").Value = "'5.741"
$ws.Range("E33").Value = '  -7.09%  '

# Row 34
$ws.Range("D34").Value = "'0.08658"
$ws.Range("E34").Value = '  -1.47%  '

# Row 35
$ws.Range("E35").Value = '  -7.00%  '

# Row 36
$ws.Range("D36").Value = "'1.519"
$ws.Range("E36").Value = '  +0.34%  '

# Row 37
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = "'5.072"
$ws.Range("E37").Value = '  -5.05%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'0.02250"
$ws.Range("E38").Value = '  -7.47%  '

# Row 39
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = "'0.6422"
$ws.Range("E39").Value = '  -6.94%  '

# Row 40
$ws.Range("D40").Value = "'0.06016"
$ws.Range("E40").Value = '  -5.38%  '

# Row 41
$ws.Range("D41").Value = "'0.2078"
$ws.Range("E41").Value = '  -5.33%  '

# Row 42
$ws.Range("E42").Value = '  -4.34%  '

# Row 43
$ws.Range("E43").Value = '  +0.85%  '

# Row 44
$ws.Range("D44").Value = "'7.901"
$ws.Range("E44").Value = '  -5.29%  '

# Row 45
$ws.Range("D45").Value = "'13.53"
$ws.Range("E45").Value = '  -5.10%  '

# Row 46
$ws.Range("D46").Value = "'3.790"
$ws.Range("E46").Value = '  -0.80%  '

# Row 47
$ws.Range("D47").Value = "'0.5945"
$ws.Range("E47").Value = '  -5.99%  '

# Row 48
$ws.Range("D48").Value = "'125.75"
$ws.Range("E48").Value = '  -4.97%  '

# Row 49
$ws.Range("D49").Value = "'1.967"
$ws.Range("E49").Value = '  -5.93%  '

# Row 50
$ws.Range("D50").Value = "'1.143"
$ws.Range("E50").Value = '  -0.39%  '

# Row 51
$ws.Range("D51").Value = "'0.06924"
$ws.Range("E51").Value = '  -6.64%  '
